$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.486268188819775
$ws.Range("C2").Value = 0.9857170294836387
$ws.Range("D2").Value = 0.5109070847733245
$ws.Range("G2").Value = 0.4760219657335256
$ws.Range("H2").Value = 0.998

# Row 3
$ws.Range("B3").Value = 0.09527669313617208
$ws.Range("C3").Value = 0.9986832503033243
$ws.Range("D3").Value = 0.2507181404971618
$ws.Range("G3").Value = 0.4760219657335256
$ws.Range("H3").Value = 0.998

# Row 4
$ws.Range("B4").Value = 0.03305850808665715
$ws.Range("C4").Value = 0.9996592730061015
$ws.Range("D4").Value = 0.1316546126141594
$ws.Range("G4").Value = 0.4760219657335256
$ws.Range("H4").Value = 0.998

# Row 5
$ws.Range("B5").Value = 0.1033622337665045
$ws.Range("C5").Value = 0.9993843006650475
$ws.Range("D5").Value = 0.2207132792383323
$ws.Range("G5").Value = 0.4760219657335256
$ws.Range("H5").Value = 0.998
